# Update quantity column (C) values for the line items, and the
# corresponding "Upto date Amount" (G/H) text values that depend on them,
# plus the two summary rows (19 and 21) that roll up the grand total.
#
# The G/H "amount" cells hold their numbers as text (e.g. "3584.00"), so
# the cell is pre-formatted as Text ("@") before the value is written —
# otherwise Excel would auto-convert the numeric-looking string back into
# a real number and drop the trailing ".00".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Qty executed upto date (no dependent amount cell - F8 is 0)
$ws.Range("C8").Value = 99

# Row 9 - Short point
$ws.Range("C9").Value = 14
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "3584.00"

# Row 10 - Medium point
$ws.Range("C10").Value = 17
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "8024.00"

# Row 11 - Long point
$ws.Range("C11").Value = 26
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "17212.00"

# Row 12 (no dependent amount cell - F12 is 0)
$ws.Range("C12").Value = 36

# Row 13 - On board
$ws.Range("C13").Value = 91
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "12376.00"

# Row 14 - P & F switch
$ws.Range("C14").Value = 86
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "1978.00"

# Row 15 - Total row (no dependent amount cell change)
$ws.Range("C15").Value = 68

# Row 16 - Add Tender Premium row (no dependent amount cell change)
$ws.Range("C16").Value = 53

# Row 17 - Grand Total row (no dependent amount cell change)
$ws.Range("C17").Value = 42

# Row 19 - Grand Total Rs. summary
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "43174.00"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "43174.00"

# Row 21 - NET PAYABLE AMOUNT Rs. summary
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "43174.00"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "43174.00"
